$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set the "Runmode" column (C2:C7) to "Y" for every test case
$ws.Range("C2:C7").Value = "Y"

# Reflect the selection left behind after filling the range
$ws.Range("C2:C7").Select()
$ws.Application.ActiveCell = $ws.Range("C2")
